$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain plain text (they contain
# values like "61.081.74" or "0.126" that Excel would otherwise
# auto-convert to numbers) while we overwrite the changed cells.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "61.081.74"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "3.394.90"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "573.11"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "142.22"
$ws.Range("E6").Value = "  -1.36%  "
$ws.Range("D7").Value = "3.396.08"
$ws.Range("E7").Value = "  -1.55%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("D13").Value = "3.973.23"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "0.126"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "27.83"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("D17").Value = "3.394.34"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "61.115.32"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("E19").Value = "  -3.99%  "
$ws.Range("E20").Value = "  -4.43%  "
$ws.Range("E21").Value = "  -5.19%  "
$ws.Range("D22").Value = "382.57"
$ws.Range("E22").Value = "  -4.15%  "
$ws.Range("D23").Value = "0.556"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").Value = "74.60"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("E26").Value = "  -4.24%  "
$ws.Range("D27").Value = "3.531.45"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "7.35"
$ws.Range("E30").Value = "  -3.85%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "7.99"
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "2.16"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("E33").Value = "  -5.49%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "23.45"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("D36").Value = "6.99"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").Value = "166.46"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("D38").Value = "3.426.86"
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("D40").Value = "1.49"
$ws.Range("E40").Value = "  -4.36%  "
$ws.Range("D41").Value = "0.0771"
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("D42").Value = "27.26"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "41.87"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").Value = "4.42"
$ws.Range("E46").Value = "  -2.76%  "
$ws.Range("E47").Value = "  -3.50%  "
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "2.473.40"
$ws.Range("E49").Value = "  -5.82%  "
$ws.Range("E50").Value = "  -2.42%  "
$ws.Range("D51").Value = "22.91"
$ws.Range("E51").Value = "  -0.76%  "

# Restore the default (unstyled) cell style now that the text values
# are safely stored, so no stray number-format styling remains.
$textRange.Style = "Normal"

